$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Val($addr, $val) {
    $ws.Range($addr).Value2 = $val
}
function Set-DateText($addr, $val) {
    # Force text storage so date-like strings ("YYYY-MM-DD") are not
    # auto-converted to date serial numbers, matching the source's
    # plain-text (inlineStr) encoding of these fields.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value2 = $val
}

# Insert two new blank rows before the old row 2 (shifts old rows 2-21 down to rows 4-23)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Row 2
Set-Val "A2" 112275278
Set-Val "B2" 90823
Set-Val "C2" "Ovaliderad"
Set-Val "D2" "NT"
Set-Val "E2" 5966
Set-Val "F2" "Motaggsvamp"
Set-Val "G2" "Sarcodon squamosus"
Set-Val "H2" "(Schaeff.) Quél."
Set-Val "P2" "Naturskoleskogen, Srm"
Set-Val "Q2" 646126
Set-Val "R2" 6568636
Set-Val "S2" 25
Set-Val "T2" "Stockholm"
Set-Val "U2" "Södertälje"
Set-Val "V2" "Södermanland"
Set-Val "W2" "Södertälje"
Set-DateText "Y2" "2023-09-23"
Set-DateText "AA2" "2023-09-23"
Set-Val "AD2" $false
Set-Val "AE2" $false
Set-Val "AG2" $false
Set-Val "AI2" "äldre barrskog"
Set-Val "AW2" "Hans Rydberg"
Set-Val "AX2" "Hans Rydberg"

# Row 3
Set-Val "A3" 112275279
Set-Val "B3" 90818
Set-Val "C3" "Ovaliderad"
Set-Val "D3" "NT"
Set-Val "E3" 4368
Set-Val "F3" "Dofttaggsvamp"
Set-Val "G3" "Hydnellum suaveolens"
Set-Val "H3" "(Scop.:Fr.) P. Karst."
Set-Val "P3" "Naturskoleskogen, Srm"
Set-Val "Q3" 646126
Set-Val "R3" 6568636
Set-Val "S3" 25
Set-Val "T3" "Stockholm"
Set-Val "U3" "Södertälje"
Set-Val "V3" "Södermanland"
Set-Val "W3" "Södertälje"
Set-DateText "Y3" "2023-09-23"
Set-DateText "AA3" "2023-09-23"
Set-Val "AD3" $false
Set-Val "AE3" $false
Set-Val "AG3" $false
Set-Val "AI3" "äldre barrskog"
Set-Val "AW3" "Hans Rydberg"
Set-Val "AX3" "Hans Rydberg"

# Append two new rows at the end (rows 24 and 25)
# Row 24
Set-Val "A24" 112275250
Set-Val "B24" 89735
Set-Val "C24" "Ovaliderad"
Set-Val "D24" "VU"
Set-Val "E24" 1106
Set-Val "F24" "Vågticka"
Set-Val "G24" "Osteina undosa"
Set-Val "H24" "(Peck) Zmitr."
Set-Val "P24" "Naturskoleskogen, Srm"
Set-Val "Q24" 646166
Set-Val "R24" 6568529
Set-Val "S24" 25
Set-Val "T24" "Stockholm"
Set-Val "U24" "Södertälje"
Set-Val "V24" "Södermanland"
Set-Val "W24" "Södertälje"
Set-DateText "Y24" "2023-09-23"
Set-DateText "AA24" "2023-09-23"
Set-Val "AD24" $false
Set-Val "AE24" $false
Set-Val "AG24" $false
Set-Val "AI24" "äldre barrskog"
Set-Val "AL24" "gran"
Set-Val "AO24" "gran"
Set-Val "AW24" "Hans Rydberg"
Set-Val "AX24" "Hans Rydberg"

# Row 25
Set-Val "A25" 112275263
Set-Val "B25" 90843
Set-Val "C25" "Ovaliderad"
Set-Val "D25" "NT"
Set-Val "E25" 5448
Set-Val "F25" "Svartvit taggsvamp"
Set-Val "G25" "Phellodon connatus"
Set-Val "H25" "(Schultz) nom.prov"
Set-Val "P25" "Naturskoleskogen, Srm"
Set-Val "Q25" 646200
Set-Val "R25" 6568598
Set-Val "S25" 50
Set-Val "T25" "Stockholm"
Set-Val "U25" "Södertälje"
Set-Val "V25" "Södermanland"
Set-Val "W25" "Södertälje"
Set-DateText "Y25" "2023-09-23"
Set-DateText "AA25" "2023-09-23"
Set-Val "AD25" $false
Set-Val "AE25" $false
Set-Val "AG25" $false
Set-Val "AI25" "äldre barrskog"
Set-Val "AW25" "Hans Rydberg"
Set-Val "AX25" "Hans Rydberg"
